$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.862.34'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.66%  '

$ws.Range("E3").Value = '  +1.26%  '

$ws.Range("E4").Value = '  +0.72%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.05%  '

$ws.Range("E6").Value = '  +0.47%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3892'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.17%  '

$ws.Range("E8").Value = '  -0.80%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.39'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.345'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.46%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.003'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.71%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08426'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.88%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.83'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.61%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.067'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.50%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.884'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.07%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001315'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.96%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.649.95'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.90%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.34%  '

$ws.Range("E19").Value = '  +0.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.83%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.925'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.93%  '

$ws.Range("E22").Value = '  +0.54%  '

$ws.Range("E23").Value = '  +1.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.860.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.68%  '

$ws.Range("E25").Value = '  -0.94%  '

$ws.Range("E26").Value = '  +3.68%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.96'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.60%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.79'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.84%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.386'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.40%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '138.28'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.73%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.772'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.84%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.511'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.32%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.831.59'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.95%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.045'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.88%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08019'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.87%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02951'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.27%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.701'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.84'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.68%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2673'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.37%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09089'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.92%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7548'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.34%  '

$ws.Range("E42").Value = '  -2.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.418'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.26%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.37'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.38%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6938'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.453'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.29%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.082'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08269'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.96'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.83%  '

$ws.Range("E51").Value = '  -0.79%  '

